$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update term labels (D column) to reflect season releveling (Winter dropped, Spring added)
$ws.Range("D5").Value = "seasonSpring"
$ws.Range("D6").Value = "seasonSummer"
$ws.Range("D7").Value = "seasonFall"
$ws.Range("D8").Value = "fish_basinWest:seasonSpring"
$ws.Range("D9").Value = "fish_basinNorth:seasonSpring"
$ws.Range("D10").Value = "fish_basinWest:seasonSummer"
$ws.Range("D11").Value = "fish_basinNorth:seasonSummer"
$ws.Range("D12").Value = "fish_basinWest:seasonFall"
$ws.Range("D13").Value = "fish_basinNorth:seasonFall"

# Update numeric estimate/std.error/statistic/p.value columns
$ws.Range("E2").Value = 0.890273349211537
$ws.Range("F2").Value = 0.0368068944451756
$ws.Range("G2").Value = 24.1876790376219
$ws.Range("H2").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000299885823158512
$ws.Range("E3").Value = 0.0855718387729713
$ws.Range("F3").Value = 0.0462083381610194
$ws.Range("G3").Value = 1.85187007753415
$ws.Range("H3").Value = 0.0640444816511403
$ws.Range("E4").Value = 0.326839706197089
$ws.Range("F4").Value = 0.0568589075625431
$ws.Range("G4").Value = 5.74825863190521
$ws.Range("H4").Value = 0.00000000901672445966544
$ws.Range("E5").Value = 0.493774947317767
$ws.Range("F5").Value = 0.0137327213317911
$ws.Range("G5").Value = 35.9560887742391
$ws.Range("H5").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000406559277235843
$ws.Range("E6").Value = 1.20658634091772
$ws.Range("F6").Value = 0.013593949360044
$ws.Range("G6").Value = 88.7590728022113
$ws.Range("H6").Value = 0
$ws.Range("E7").Value = 1.1279933706489
$ws.Range("F7").Value = 0.0132717258152763
$ws.Range("G7").Value = 84.992214754055
$ws.Range("H7").Value = 0
$ws.Range("E8").Value = -0.021232815987736
$ws.Range("F8").Value = 0.0186517741610287
$ws.Range("G8").Value = -1.13838049959345
$ws.Range("H8").Value = 0.254961633149316
$ws.Range("E9").Value = -0.273984351683324
$ws.Range("F9").Value = 0.0218687244996009
$ws.Range("G9").Value = -12.5285931371226
$ws.Range("H9").Value = 0.00000000000000000000000000000000000520782920806283
$ws.Range("E10").Value = -0.0468302195593423
$ws.Range("F10").Value = 0.0183477409328543
$ws.Range("G10").Value = -2.55236978387274
$ws.Range("H10").Value = 0.0106992895461227
$ws.Range("E11").Value = -0.450130254304272
$ws.Range("F11").Value = 0.0224021836819427
$ws.Range("G11").Value = -20.0931418425561
$ws.Range("H11").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000847272504723124
$ws.Range("E12").Value = -0.0822485260186274
$ws.Range("F12").Value = 0.0180131710045114
$ws.Range("G12").Value = -4.5660214960502
$ws.Range("H12").Value = 0.00000497067818095543
$ws.Range("E13").Value = -0.388174676513601
$ws.Range("F13").Value = 0.0214361986366336
$ws.Range("G13").Value = -18.1083728087044
$ws.Range("H13").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000273720259666219
$ws.Range("E14").Value = 0.0867992057914639
$ws.Range("E15").Value = 0.0255597934641572
